# Rename the "Cross references" sheet to "Database references" (in every
# model this pattern occurs) because the old name is ambiguous, and make
# that sheet the active/selected tab, as it was when the workbook was saved.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Cross references") {
        $ws.Name = "Database references"
        $ws.Activate()
        $ws.Select()
    }
}
